# edit.ps1
# Applies the "Updated symbol list" GitHub Actions commit to cryptos.xlsx.
# For every data row (2-51) the scraped Price (D), Volume(1h) (E) and Hora
# (G) columns are refreshed. All three columns hold plain text in the
# source workbook (numeric-looking strings, percentages, hour-of-day), so
# each cell's NumberFormat is forced to "@" (Text) before the value is
# written — this stops Excel's COM layer from auto-coercing strings such
# as "308.91" or "0.23%" into numbers/percentages (which would also
# introduce floating-point noise like 308.91000000000003).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@(2, '308.91', '0.23%', '20')
    ,@(3, '40.90', '-0.09%', '20')
    ,@(4, '5.126', '1.41%', '20')
    ,@(5, '0.07617', '-0.38%', '20')
    ,@(6, '1.605', '-0.49%', '20')
    ,@(7, '2.484', '3.49%', '20')
    ,@(8, '0.9079', $null, '20')
    ,@(9, '0.1254', '24.50%', '20')
    ,@(10, '0.1806', '2.70%', '20')
    ,@(11, '0.09093', '0.28%', '20')
    ,@(12, '0.04297', '-2.37%', '20')
    ,@(13, $null, '-0.86%', '20')
    ,@(14, '0.001257', $null, '20')
    ,@(15, '0.005725', '-1.50%', '20')
    ,@(16, $null, $null, '20')
    ,@(17, '4.282', '0.56%', '20')
    ,@(18, $null, '1.43%', '20')
    ,@(19, '6.940', '2.89%', '20')
    ,@(20, '0.1393', '2.70%', '20')
    ,@(21, '0.2688', '-1.22%', '20')
    ,@(22, '0.04035', '-2.84%', '20')
    ,@(23, '0.001268', '3.99%', '20')
    ,@(24, '0.004046', '-0.85%', '20')
    ,@(25, $null, '-1.99%', '20')
    ,@(26, $null, '24.83%', '20')
    ,@(27, $null, $null, '20')
    ,@(28, $null, $null, '20')
    ,@(29, $null, $null, '20')
    ,@(30, $null, $null, '20')
    ,@(31, $null, $null, '20')
    ,@(32, $null, $null, '20')
    ,@(33, $null, $null, '20')
    ,@(34, $null, $null, '20')
    ,@(35, $null, $null, '20')
    ,@(36, $null, $null, '20')
    ,@(37, $null, $null, '20')
    ,@(38, '0.02419', $null, '20')
    ,@(39, '0.05229', '0.95%', '20')
    ,@(40, '0.007845', '0.51%', '20')
    ,@(41, '0.1302', '-0.82%', '20')
    ,@(42, '0.006798', '-3.99%', '20')
    ,@(43, '0.001903', '-2.17%', '20')
    ,@(44, '0.007449', '-7.03%', '20')
    ,@(45, '0.3362', '9.87%', '20')
    ,@(46, '0.00006899', '8.34%', '20')
    ,@(47, $null, '0.24%', '20')
    ,@(48, '0.1301', '2,306.37%', '20')
    ,@(49, $null, '-31.71%', '20')
    ,@(50, '0.00002104', '0.24%', '20')
    ,@(51, $null, '0.24%', '20')
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $price  = $r[1]
    $volume = $r[2]
    $hour   = $r[3]

    if ($null -ne $price) {
        $cell = $ws.Cells.Item($rowNum, 4)   # column D - Price
        $cell.NumberFormat = "@"
        $cell.Value = $price
    }

    if ($null -ne $volume) {
        $cell = $ws.Cells.Item($rowNum, 5)   # column E - Volume(1h)
        $cell.NumberFormat = "@"
        $cell.Value = $volume
    }

    if ($null -ne $hour) {
        $cell = $ws.Cells.Item($rowNum, 7)   # column G - Hora
        $cell.NumberFormat = "@"
        $cell.Value = $hour
    }
}
